# The sheet currently holds a 2-column table ("price"/"demand") anchored at
# B2:C9. The edit re-anchors that same table one row up and one column left,
# so it lands at A1:B8 instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row + all data rows, in final (target) layout.
$table = @(
    @("price", "demand"),
    @(5,  600),
    @(10, 550),
    @(15, 700),
    @(17, 680),
    @(20, 500),
    @(23, 400),
    @(25, 250)
)

# Remove the old table (B2:C9) now that its contents have been captured above.
$ws.Range("B2:C9").Clear()

# Write the table back out starting at A1.
for ($i = 0; $i -lt $table.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $table[$i][0]
    $ws.Cells.Item($row, 2).Value = $table[$i][1]
}

# Collapse the selection back down to a single cell (A1) instead of the
# stale "B1:E1048576" selection that was left over from the prior edit.
$ws.Range("A1").Select()
